# [AA | 13/4/2018] : commit for Sample tab changes and profile page
#
# Renames the original (and only) sheet "Sheet1" to "Sample_Data", then adds
# three new sheets: "Fields_Data" (field metadata), "Profile_Data" (per-field
# profiling summary) and "Categorical_Data" (frequency breakdown for one
# categorical field). Final tab order is Sample_Data, Profile_Data,
# Categorical_Data, Fields_Data, with Categorical_Data left as the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the existing sheet
# ---------------------------------------------------------------------------
$sampleData = $wb.Worksheets.Item(1)
$sampleData.Name = "Sample_Data"

# ---------------------------------------------------------------------------
# 2) Add "Fields_Data" right after Sample_Data (picks up sheetId 2)
# ---------------------------------------------------------------------------
$fieldsData = $wb.Worksheets.Add($null, $sampleData)
$fieldsData.Name = "Fields_Data"

$fieldsData.Range("A1").Value = 1
$fieldsData.Range("B1").Value = "ID"
$fieldsData.Range("C1").Value = "INTEGER"
$fieldsData.Range("F1").Value = "No"
$fieldsData.Range("G1").Value = "No"

$fieldsData.Range("A2").Value = 2
$fieldsData.Range("B2").Value = "NAME"
$fieldsData.Range("C2").Value = "STRING"
$fieldsData.Range("F2").Value = "No"
$fieldsData.Range("G2").Value = "No"

$fieldsData.Range("A3").Value = 3
$fieldsData.Range("B3").Value = "AGE"
$fieldsData.Range("C3").Value = "INTEGER"
$fieldsData.Range("F3").Value = "No"
$fieldsData.Range("G3").Value = "No"

$fieldsData.Range("A4").Value = 4
$fieldsData.Range("B4").Value = "ADDRESS"
$fieldsData.Range("C4").Value = "STRING"
$fieldsData.Range("F4").Value = "No"
$fieldsData.Range("G4").Value = "No"

$fieldsData.Range("A5").Value = 5
$fieldsData.Range("B5").Value = "SALARY"
$fieldsData.Range("C5").Value = "DOUBLE"
$fieldsData.Range("F5").Value = "No"
$fieldsData.Range("G5").Value = "No"

$fieldsData.Range("F5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Add "Profile_Data" right after Fields_Data (picks up sheetId 3)
# ---------------------------------------------------------------------------
$profileData = $wb.Worksheets.Add($null, $fieldsData)
$profileData.Name = "Profile_Data"

$profileData.Range("A1").Value = "Field Name"
$profileData.Range("B1").Value = "Data Type"
$profileData.Range("C1").Value = "Scale Type"
$profileData.Range("D1").Value = "Distinct Count"
$profileData.Range("E1").Value = "Unique Count(%)"
$profileData.Range("F1").Value = "Missing(%)"

$profileData.Range("A2").Value = "id"
$profileData.Range("B2").Value = "INT"
$profileData.Range("C2").Value = "categorical"
$profileData.Range("D2").Value = 9
$profileData.Range("E2").Value = 90
$profileData.Range("F2").Value = 10

$profileData.Range("A3").Value = "name"
$profileData.Range("B3").Value = "STRING"
$profileData.Range("C3").Value = "categorical"
$profileData.Range("D3").Value = 10
$profileData.Range("E3").Value = 100
$profileData.Range("F3").Value = 0

$profileData.Range("A4").Value = "age"
$profileData.Range("B4").Value = "INT"
$profileData.Range("C4").Value = "categorical"
$profileData.Range("D4").Value = 6
$profileData.Range("E4").Value = 60
$profileData.Range("F4").Value = 10

$profileData.Range("A5").Value = "address"
$profileData.Range("B5").Value = "STRING"
$profileData.Range("C5").Value = "categorical"
$profileData.Range("D5").Value = 8
$profileData.Range("E5").Value = 80
$profileData.Range("F5").Value = 0

$profileData.Range("A6").Value = "salary"
$profileData.Range("B6").Value = "DOUBLE"
$profileData.Range("C6").Value = "categorical"
$profileData.Range("D6").Value = 8
$profileData.Range("E6").Value = 80
$profileData.Range("F6").Value = 10

$profileData.Range("A7").Value = "jobinstanceid"
$profileData.Range("B7").Value = "STRING"
$profileData.Range("C7").Value = "constant"
$profileData.Range("D7").Value = 1
$profileData.Range("E7").Value = 10
$profileData.Range("F7").Value = 0

$profileData.Columns.Item(1).ColumnWidth = 13.14
$profileData.Columns.Item(2).ColumnWidth = 9.71
$profileData.Columns.Item(3).ColumnWidth = 10.57
$profileData.Columns.Item(4).ColumnWidth = 13.57
$profileData.Columns.Item(5).ColumnWidth = 16.29
$profileData.Columns.Item(6).ColumnWidth = 10.71

$profileData.Range("D7").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4) Add "Categorical_Data" right after Profile_Data (picks up sheetId 4)
# ---------------------------------------------------------------------------
$categoricalData = $wb.Worksheets.Add($null, $profileData)
$categoricalData.Name = "Categorical_Data"

$categoricalData.Range("A1").Value = "Data"
$categoricalData.Range("B1").Value = "Frequency"
$categoricalData.Range("C1").Value = "Percentage"

$catRows = @(
    @(777, 1, 0.1111),
    @(666, 1, 0.1111),
    @(555, 1, 0.1111),
    @(500, 1, 0.1111),
    @(444, 1, 0.1111),
    @(333, 1, 0.1111),
    @(222, 1, 0.1111),
    @(200, 1, 0.1111),
    @(100, 1, 0.1111)
)

$r = 2
foreach ($row in $catRows) {
    $categoricalData.Range("A$r").Value = $row[0]
    $categoricalData.Range("B$r").Value = $row[1]
    $categoricalData.Range("C$r").Value = $row[2]
    $categoricalData.Range("C$r").NumberFormat = "0.00%"
    $r = $r + 1
}

$categoricalData.Columns.Item(2).ColumnWidth = 10.29
$categoricalData.Columns.Item(3).ColumnWidth = 11

$categoricalData.Range("G7").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5) Move "Fields_Data" to the end of the tab strip
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$fieldsData.Move($null, $lastSheet)

# ---------------------------------------------------------------------------
# 6) Make "Categorical_Data" the active tab
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Categorical_Data").Activate()
